$d = $word.ActiveDocument

# 1) "...questions with this pr." -> "...questions with the specific product."
$d.Content.Find.Execute(
    " and provide one-to-one online chat if customers have some questions with this pr",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " and provide one-to-one online chat if customers have some questions with the specific product",
    2) | Out-Null

# 2) Remove the two fully-blank paragraphs that sit between the "Moreover..." paragraph
#    and the "According to statistics..." paragraph, and append a new closing sentence
#    to the "Moreover..." paragraph.
$moreover = $d.Content.Find
$found = $d.Content.Find.Execute("Moreover, after the order is completed, the user can comment on the product and upload pictures, which can help other users to choose the product.")

# Locate the paragraph that contains that sentence and append the new sentence to it.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*help other users to choose the product.*") {
        $p.Range.InsertAfter(" In the following, we will illustrate the reasons why we choose mobile e-commerce applications.")
        break
    }
}

# Now delete the two (now three, but only two are truly empty) blank paragraphs that
# immediately follow it.
$idx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text.Trim()
    if ($t -eq "" -and $d.Paragraphs($i).Range.Text -ne [char]13) {
        # placeholder, not used
    }
}

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*illustrate the reasons why we choose mobile e-commerce applications.*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $next1 = $target.Next()
    $next1.Range.Delete() | Out-Null
    $next2 = $target.Next()
    $next2.Range.Delete() | Out-Null
}
